$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Flow Group")

# New full list of rows (A = Name, B = Description) for rows 2..20, sorted alphabetically
# with the renamed "Emission: Total" -> "Emission: Total Rh" and several new flow-group rows added.
$rows = @(
    @{ A = "Biomass Turnover: Total"; B = "" },
    @{ A = "Decay: Total"; B = "" },
    @{ A = "Emission: Total Rh"; B = "Total carbon emissions from all ecosystem components (DOM+ Biomass)" },
    @{ A = "LULC: Emission"; B = "" },
    @{ A = "LULC: Emission CH4"; B = "" },
    @{ A = "LULC: Emission CO"; B = "" },
    @{ A = "LULC: Emission CO2"; B = "" },
    @{ A = "LULC: Emission DOM"; B = "" },
    @{ A = "LULC: Emission Live"; B = "" },
    @{ A = "LULC: Harvest"; B = "" },
    @{ A = "LULC: Mortality"; B = "" },
    @{ A = "LULC: Transfer"; B = "" },
    @{ A = "Net Biome Productivity"; B = "" },
    @{ A = "Net Ecosystem Productivity"; B = "" },
    @{ A = "Net Growth: Total"; B = "Net biomass increment before losses from disturbances" },
    @{ A = "Net Primary Productivity"; B = "" },
    @{ A = "Q10 Fast Flows"; B = "" },
    @{ A = "Q10 Slow Flows"; B = "" },
    @{ A = "Transfer: Total"; B = "" }
)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $a = $rows[$i].A
    $b = $rows[$i].B

    $ws.Cells.Item($r, 1).Value2 = $a
    if ($b -ne "") {
        $ws.Cells.Item($r, 2).Value2 = $b
    } else {
        $ws.Cells.Item($r, 2).Value2 = "placeholder"
        $ws.Cells.Item($r, 2).ClearContents()
    }
}

# Reset formatting on the populated rows (A2:B20) back to the default "Normal" style,
# matching the target, which no longer carries the bold/indented style on these rows.
$ws.Range("A2:B20").Style = "Normal"

# Un-hide columns C and D (they keep their widths but are no longer hidden).
$ws.Columns.Item(3).Hidden = $false
$ws.Columns.Item(4).Hidden = $false

# Column A needs to widen to fit the longer text now in it (AutoFit mirrors Excel's bestFit behavior).
$ws.Columns.Item(1).AutoFit() | Out-Null

# Move the active selection, matching the end-user's last edit location.
$ws.Range("B19").Select() | Out-Null
